$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Player")
$ws1.Cells.Item(1,1).Value = "캐릭터ID"
$ws1.Cells.Item(1,2).Value = "캐릭터 모델"
$ws1.Cells.Item(1,3).Value = "체력"
$ws1.Cells.Item(1,4).Value = "공격력"
$ws1.Cells.Item(1,5).Value = "방어력"
$ws1.Cells.Item(2,1).Value = 1
$ws1.Cells.Item(2,2).Value = 102
$ws1.Cells.Item(2,3).Value = 5000
$ws1.Cells.Item(2,4).Value = 500
$ws1.Cells.Item(2,5).Value = 80
$ws1.Cells.Item(3,1).Value = 2
$ws1.Cells.Item(3,2).Value = 101
$ws1.Cells.Item(3,3).Value = 3000
$ws1.Cells.Item(3,4).Value = 400
$ws1.Cells.Item(3,5).Value = 70
$ws1.Cells.Item(4,1).Value = 3
$ws1.Cells.Item(4,2).Value = 103
$ws1.Cells.Item(4,3).Value = 3000
$ws1.Cells.Item(4,4).Value = 450
$ws1.Cells.Item(4,5).Value = 60

$ws2 = $wb.Worksheets.Item("Attack")
$ws2.Cells.Item(1,1).Value = "스킬ID"
$ws2.Cells.Item(1,2).Value = "스킬이름"
$ws2.Cells.Item(1,3).Value = "스킬 타입(eAttackType)"
$ws2.Cells.Item(1,4).Value = "스킬 추가 파워"
$ws2.Cells.Item(1,5).Value = "시전가능거리"
$ws2.Cells.Item(1,6).Value = "쿨타임"
$ws2.Cells.Item(1,7).Value = "이펙트 프리팹 이름"
$ws2.Cells.Item(1,8).Value = "UI 이미지 이름"
$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "EMPTY"
$ws2.Cells.Item(2,3).Value = 0
$ws2.Cells.Item(2,4).Value = 0
$ws2.Cells.Item(2,5).Value = 0
$ws2.Cells.Item(2,6).Value = 0
$ws2.Cells.Item(2,7).Value = "NULL"
$ws2.Cells.Item(3,1).Value = 10001
$ws2.Cells.Item(3,2).Value = "기본공격(풀)"
$ws2.Cells.Item(3,3).Value = 1
$ws2.Cells.Item(3,4).Value = 1
$ws2.Cells.Item(3,5).Value = 5
$ws2.Cells.Item(3,6).Value = 3
$ws2.Cells.Item(3,7).Value = "TinyGrassEffect"
$ws2.Cells.Item(4,1).Value = 10002
$ws2.Cells.Item(4,2).Value = "기본공격(불)"
$ws2.Cells.Item(4,3).Value = 1
$ws2.Cells.Item(4,4).Value = 1
$ws2.Cells.Item(4,5).Value = 5
$ws2.Cells.Item(4,6).Value = 3
$ws2.Cells.Item(4,7).Value = "TinyFireEffect"
$ws2.Cells.Item(5,1).Value = 10003
$ws2.Cells.Item(5,2).Value = "기본공격(물)"
$ws2.Cells.Item(5,3).Value = 1
$ws2.Cells.Item(5,4).Value = 1
$ws2.Cells.Item(5,5).Value = 5
$ws2.Cells.Item(5,6).Value = 3
$ws2.Cells.Item(5,7).Value = "TinyWaterEffect"
$ws2.Cells.Item(6,1).Value = 10101
$ws2.Cells.Item(6,2).Value = "잎날가르기"
$ws2.Cells.Item(6,3).Value = 2
$ws2.Cells.Item(6,4).Value = 1.5
$ws2.Cells.Item(6,5).Value = 10
$ws2.Cells.Item(6,6).Value = 5
$ws2.Cells.Item(6,7).Value = "BigGrassEffect"
$ws2.Cells.Item(6,8).Value = "Skill_GrassThrow"
$ws2.Cells.Item(7,1).Value = 10102
$ws2.Cells.Item(7,2).Value = "불꽃세례"
$ws2.Cells.Item(7,3).Value = 2
$ws2.Cells.Item(7,4).Value = 1.5
$ws2.Cells.Item(7,5).Value = 10
$ws2.Cells.Item(7,6).Value = 5
$ws2.Cells.Item(7,7).Value = "BigFireEffect"
$ws2.Cells.Item(7,8).Value = "Skill_FireThrow"
$ws2.Cells.Item(8,1).Value = 10103
$ws2.Cells.Item(8,2).Value = "물대포"
$ws2.Cells.Item(8,3).Value = 2
$ws2.Cells.Item(8,4).Value = 1.5
$ws2.Cells.Item(8,5).Value = 10
$ws2.Cells.Item(8,6).Value = 5
$ws2.Cells.Item(8,7).Value = "BigWaterEffect"
$ws2.Cells.Item(8,8).Value = "Skill_WaterThrow"
$ws2.Cells.Item(9,1).Value = 10202
$ws2.Cells.Item(9,2).Value = "회오리불꽃"
$ws2.Cells.Item(9,3).Value = 3
$ws2.Cells.Item(9,4).Value = 1.7
$ws2.Cells.Item(9,5).Value = 3
$ws2.Cells.Item(9,6).Value = 5
$ws2.Cells.Item(9,7).Value = "BigFireEffect"
$ws2.Cells.Item(9,8).Value = "Skill_FireSpin"

$ws3 = $wb.Worksheets.Item("Model")
$ws3.Cells.Item(1,1).Value = "모델ID"
$ws3.Cells.Item(1,2).Value = "캐릭터 이름"
$ws3.Cells.Item(1,3).Value = "보유스킬 "
$ws3.Cells.Item(1,4).Value = "모델이름"
$ws3.Cells.Item(1,5).Value = "모델 UI 이미지"
$ws3.Cells.Item(2,1).Value = 101
$ws3.Cells.Item(2,2).Value = "초록이"
$ws3.Cells.Item(2,3).Value = "10001/10101"
$ws3.Cells.Item(2,4).Value = "Flora"
$ws3.Cells.Item(2,5).Value = "Model_Flora"
$ws3.Cells.Item(3,1).Value = 102
$ws3.Cells.Item(3,2).Value = "빨강이"
$ws3.Cells.Item(3,3).Value = "10002/10102"
$ws3.Cells.Item(3,4).Value = "Fiery"
$ws3.Cells.Item(3,5).Value = "Model_Fiery"
$ws3.Cells.Item(4,1).Value = 103
$ws3.Cells.Item(4,2).Value = "파랑이"
$ws3.Cells.Item(4,3).Value = "10003/10103"
$ws3.Cells.Item(4,4).Value = "Cyclopes"
$ws3.Cells.Item(4,5).Value = "Model_Cyclopes"
$ws3.Cells.Item(5,1).Value = 104
$ws3.Cells.Item(5,3).Value = 10001
$ws3.Cells.Item(5,4).Value = "Beezee"
$ws3.Cells.Item(6,1).Value = 105
$ws3.Cells.Item(6,3).Value = 10001
$ws3.Cells.Item(6,4).Value = "Kaktos"
$ws3.Cells.Item(7,1).Value = 106
$ws3.Cells.Item(7,3).Value = 10001
$ws3.Cells.Item(7,4).Value = "Ketchup"
$ws3.Cells.Item(8,1).Value = 107
$ws3.Cells.Item(8,3).Value = 10001
$ws3.Cells.Item(8,4).Value = "Woody"

$ws4 = $wb.Worksheets.Item("Enemy")
$ws4.Cells.Item(1,1).Value = "맵 레벨"
$ws4.Cells.Item(1,2).Value = "출현 모델"
$ws4.Cells.Item(1,3).Value = "적 체력"
$ws4.Cells.Item(1,4).Value = "적 공격력"
$ws4.Cells.Item(1,5).Value = "적 방어력"
$ws4.Cells.Item(2,1).Value = 100
$ws4.Cells.Item(2,2).Value = "104/105/106/107"
$ws4.Cells.Item(2,3).Value = 500
$ws4.Cells.Item(2,4).Value = 100
$ws4.Cells.Item(2,5).Value = 30

# Column width adjustments (Attack sheet: new G/H columns for prefab / UI image name)
$ws2.Columns.Item(5).ColumnWidth = 11.285714285714286
$ws2.Columns.Item(6).ColumnWidth = 11.285714285714286
$ws2.Columns.Item(7).ColumnWidth = 15.684151785714286
$ws2.Columns.Item(8).ColumnWidth = 14.184151785714286

# Column width adjustments (Model sheet: new E column for model UI image)
$ws3.Columns.Item(4).ColumnWidth = 10.383370535714286
$ws3.Columns.Item(5).ColumnWidth = 19.484933035714285

# Page setup (Attack sheet)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selections
$ws1.Range("F11").Select()
$ws2.Range("G10").Select()
$ws3.Range("C10").Select()
$ws4.Activate()
$ws4.Range("E13").Select()
